$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.200.37"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "2.422.16"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.24"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.72"
$ws.Range("E6").Value = "  +2.89%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.535"
$ws.Range("E8").Value = "  +1.48%  "
$ws.Range("D9").Value = "2.419.35"
$ws.Range("E9").Value = "  +1.22%  "
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("E11").Value = "  -1.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.42"
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.10"
$ws.Range("E14").Value = "  +2.23%  "
$ws.Range("E15").Value = "  +5.12%  "
$ws.Range("D16").Value = "2.859.55"
$ws.Range("E16").Value = "  +2.05%  "
$ws.Range("D17").Value = "62.067.13"
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("D18").Value = "2.420.03"
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.27"
$ws.Range("E19").Value = "  +2.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "325.67"
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.78"
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.43"
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("E25").Value = "  -0.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.07"
$ws.Range("E26").Value = "  +3.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "589.18"
$ws.Range("E27").Value = "  +12.15%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.539.47"
$ws.Range("E28").Value = "  +1.57%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("D30").Value = "0.0₃0947"
$ws.Range("E30").Value = "  +3.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.27"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.46"
$ws.Range("E32").Value = "  +4.77%  "
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("E34").Value = "  +2.34%  "
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.74"
$ws.Range("E36").Value = "  +3.44%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.80"
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "154.00"
$ws.Range("E39").Value = "  +4.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.385"
$ws.Range("E40").Value = "  +1.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.71"
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.84"
$ws.Range("E42").Value = "  -4.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.36"
$ws.Range("E44").Value = "  +8.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "150.49"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.66"
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("E47").Value = "  +2.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.47"
$ws.Range("E48").Value = "  +3.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.593"
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0927"
$ws.Range("E50").Value = "  +1.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0230"
$ws.Range("E51").Value = "  +1.93%  "
